$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 98.75
$ws.Cells.Item(12, 9).Value = 81.666664
$ws.Cells.Item(12, 10).Value = 150
$ws.Cells.Item(12, 11).Value = 81.666664
$ws.Cells.Item(12, 12).Value = 150
$ws.Cells.Item(12, 13).Value = 88.333336
$ws.Cells.Item(12, 14).Value = -490

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 265
$ws.Cells.Item(28, 9).Value = 265
$ws.Cells.Item(28, 11).Value = 265
$ws.Cells.Item(28, 13).Value = 220

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 435.25
$ws.Cells.Item(33, 9).Value = 484.75
$ws.Cells.Item(33, 10).Value = 187.75
$ws.Cells.Item(33, 11).Value = 484.75
$ws.Cells.Item(33, 12).Value = 187.75
$ws.Cells.Item(33, 13).Value = -255.75
$ws.Cells.Item(33, 14).Value = -645.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1620.4166
$ws.Cells.Item(70, 9).Value = 750
$ws.Cells.Item(70, 11).Value = 2250
$ws.Cells.Item(70, 13).Value = -1980

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 1620.4166
$ws.Cells.Item(73, 9).Value = 750
$ws.Cells.Item(73, 11).Value = 2250
$ws.Cells.Item(73, 13).Value = -1314

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 12632.5
$ws.Cells.Item(74, 9).Value = 12632.5
$ws.Cells.Item(74, 11).Value = 12632.5
$ws.Cells.Item(74, 13).Value = -11696.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 12632.5
$ws.Cells.Item(77, 9).Value = 12632.5
$ws.Cells.Item(77, 11).Value = 63162.5
$ws.Cells.Item(77, 13).Value = -58482.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 4049.3333
$ws.Cells.Item(80, 9).Value = 1649.5
$ws.Cells.Item(80, 11).Value = 4948.5
$ws.Cells.Item(80, 13).Value = -3950.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 8).Value = 4049.3333
$ws.Cells.Item(83, 9).Value = 1649.5
$ws.Cells.Item(83, 11).Value = 14845.5
$ws.Cells.Item(83, 13).Value = -9853.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(87, 8).Value = 69975
$ws.Cells.Item(87, 10).Value = 69975
$ws.Cells.Item(87, 12).Value = 69975
$ws.Cells.Item(87, 14).Value = -72471

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(90, 8).Value = 69975
$ws.Cells.Item(90, 10).Value = 69975
$ws.Cells.Item(90, 12).Value = 209925
$ws.Cells.Item(90, 14).Value = -222405

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 1536.5
$ws.Cells.Item(107, 9).Value = 1509.2354
$ws.Cells.Item(107, 11).Value = 1509.2354
$ws.Cells.Item(107, 13).Value = 410.7646

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2608
$ws.Cells.Item(137, 9).Value = 2489.25
$ws.Cells.Item(137, 11).Value = 7467.75
$ws.Cells.Item(137, 13).Value = -4917.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1937.8572
$ws.Cells.Item(2, 9).Value = 1937.8572
$ws.Cells.Item(2, 11).Value = 1937.8572
$ws.Cells.Item(2, 13).Value = -1824.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7906.5
$ws.Cells.Item(32, 9).Value = 7906.5
$ws.Cells.Item(32, 11).Value = 7906.5
$ws.Cells.Item(32, 13).Value = -7619.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(41, 8).Value = 793.3333
$ws.Cells.Item(41, 9).Value = 793.3333
$ws.Cells.Item(41, 11).Value = 793.3333
$ws.Cells.Item(41, 13).Value = -379.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(62, 8).Value = 40000
$ws.Cells.Item(62, 10).Value = 40000
$ws.Cells.Item(62, 12).Value = 40000
$ws.Cells.Item(62, 14).Value = -41248

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(65, 8).Value = 40000
$ws.Cells.Item(65, 10).Value = 40000
$ws.Cells.Item(65, 12).Value = 120000
$ws.Cells.Item(65, 14).Value = -126240

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 2464.9
$ws.Cells.Item(88, 9).Value = 1949.75
$ws.Cells.Item(88, 10).Value = 2808.3333
$ws.Cells.Item(88, 11).Value = 1949.75
$ws.Cells.Item(88, 12).Value = 2808.3333
$ws.Cells.Item(88, 13).Value = -1543.75
$ws.Cells.Item(88, 14).Value = -3620.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 2464.9
$ws.Cells.Item(91, 9).Value = 1949.75
$ws.Cells.Item(91, 10).Value = 2808.3333
$ws.Cells.Item(91, 11).Value = 1949.75
$ws.Cells.Item(91, 12).Value = 2808.3333
$ws.Cells.Item(91, 13).Value = -545.75
$ws.Cells.Item(91, 14).Value = -5616.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 3597
$ws.Cells.Item(102, 9).Value = 3963
$ws.Cells.Item(102, 10).Value = 2499
$ws.Cells.Item(102, 11).Value = 3963
$ws.Cells.Item(102, 12).Value = 2499
$ws.Cells.Item(102, 13).Value = -2341
$ws.Cells.Item(102, 14).Value = -5743

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1937.8572
$ws.Cells.Item(116, 9).Value = 1937.8572
$ws.Cells.Item(116, 11).Value = 1937.8572
$ws.Cells.Item(116, 13).Value = 356.1428000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1937.8572
$ws.Cells.Item(3, 9).Value = 1937.8572
$ws.Cells.Item(3, 11).Value = 1937.8572
$ws.Cells.Item(3, 13).Value = -1823.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(98, 8).Value = 650000
$ws.Cells.Item(98, 10).Value = 650000
$ws.Cells.Item(98, 12).Value = 650000
$ws.Cells.Item(98, 14).Value = -655990

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 21201.4
$ws.Cells.Item(99, 9).Value = 26001.75
$ws.Cells.Item(99, 11).Value = 26001.75
$ws.Cells.Item(99, 13).Value = -24503.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2145
$ws.Cells.Item(105, 9).Value = 2145
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 2145
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -398
$ws.Cells.Item(105, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 13999
$ws.Cells.Item(62, 9).Value = 19998.5
$ws.Cells.Item(62, 10).Value = 2000
$ws.Cells.Item(62, 11).Value = 19998.5
$ws.Cells.Item(62, 12).Value = 2000
$ws.Cells.Item(62, 13).Value = -19374.5
$ws.Cells.Item(62, 14).Value = -3248

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 13999
$ws.Cells.Item(65, 9).Value = 19998.5
$ws.Cells.Item(65, 10).Value = 2000
$ws.Cells.Item(65, 11).Value = 99992.5
$ws.Cells.Item(65, 12).Value = 10000
$ws.Cells.Item(65, 13).Value = -96872.5
$ws.Cells.Item(65, 14).Value = -16240

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 69996.25
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 13).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(71, 8).Value = 69996.25
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 13).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 3908.3333
$ws.Cells.Item(99, 9).Value = 3909
$ws.Cells.Item(99, 10).Value = 3907
$ws.Cells.Item(99, 11).Value = 3909
$ws.Cells.Item(99, 12).Value = 3907
$ws.Cells.Item(99, 13).Value = -2411
$ws.Cells.Item(99, 14).Value = -6903

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 3908.3333
$ws.Cells.Item(126, 9).Value = 3909
$ws.Cells.Item(126, 10).Value = 3907
$ws.Cells.Item(126, 11).Value = 11727
$ws.Cells.Item(126, 12).Value = 11721
$ws.Cells.Item(126, 13).Value = -9257
$ws.Cells.Item(126, 14).Value = -16661

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 349.5
$ws.Cells.Item(26, 9).Value = 349.5
$ws.Cells.Item(26, 11).Value = 1048.5
$ws.Cells.Item(26, 13).Value = -760.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 19994.834
$ws.Cells.Item(39, 10).Value = 19994.834
$ws.Cells.Item(39, 12).Value = 59984.50199999999
$ws.Cells.Item(39, 14).Value = -60572.50199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 22000
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 13).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 8).Value = 22000
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 13).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 2961.1428
$ws.Cells.Item(98, 9).Value = 2950.6
$ws.Cells.Item(98, 11).Value = 8851.799999999999
$ws.Cells.Item(98, 13).Value = -7353.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 1666
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 4665
$ws.Cells.Item(131, 9).Value = 3995
$ws.Cells.Item(131, 10).Value = 5000
$ws.Cells.Item(131, 11).Value = 11985
$ws.Cells.Item(131, 12).Value = 15000
$ws.Cells.Item(131, 13).Value = -6945
$ws.Cells.Item(131, 14).Value = -25080

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 30000
$ws.Cells.Item(15, 10).Value = 30000
$ws.Cells.Item(15, 12).Value = 30000
$ws.Cells.Item(15, 14).Value = -30576

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3577.6
$ws.Cells.Item(80, 10).Value = 3496
$ws.Cells.Item(80, 12).Value = 3496
$ws.Cells.Item(80, 14).Value = -5492

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(81, 8).Value = 30000
$ws.Cells.Item(81, 10).Value = 30000
$ws.Cells.Item(81, 12).Value = 30000
$ws.Cells.Item(81, 14).Value = -31996

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3577.6
$ws.Cells.Item(83, 10).Value = 3496
$ws.Cells.Item(83, 12).Value = 17480
$ws.Cells.Item(83, 14).Value = -27464

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(84, 8).Value = 30000
$ws.Cells.Item(84, 10).Value = 30000
$ws.Cells.Item(84, 12).Value = 90000
$ws.Cells.Item(84, 14).Value = -99984

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(99, 8).Value = 9997
$ws.Cells.Item(99, 9).Value = 9997
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 9997
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 14).ClearContents()
$ws.Cells.Item(99, 13).Value = -7751

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 8298.25
$ws.Cells.Item(122, 9).Value = 2000
$ws.Cells.Item(122, 11).Value = 6000
$ws.Cells.Item(122, 13).Value = -3550

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3809.8572
$ws.Cells.Item(7, 9).Value = 3809.8572
$ws.Cells.Item(7, 11).Value = 3809.8572
$ws.Cells.Item(7, 13).Value = -3697.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2128.8
$ws.Cells.Item(16, 9).Value = 2241
$ws.Cells.Item(16, 10).Value = 1867
$ws.Cells.Item(16, 11).Value = 2241
$ws.Cells.Item(16, 12).Value = 1867
$ws.Cells.Item(16, 13).Value = -2071
$ws.Cells.Item(16, 14).Value = -2207

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3000
$ws.Cells.Item(40, 10).Value = 3000
$ws.Cells.Item(40, 12).Value = 3000
$ws.Cells.Item(40, 14).Value = -3272

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 3809.8572
$ws.Cells.Item(126, 9).Value = 3809.8572
$ws.Cells.Item(126, 11).Value = 11429.5716
$ws.Cells.Item(126, 13).Value = -8959.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 5991.5
$ws.Cells.Item(132, 9).Value = 5989
$ws.Cells.Item(132, 11).Value = 17967
$ws.Cells.Item(132, 13).Value = -15437

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1194.5
$ws.Cells.Item(122, 9).Value = 1189
$ws.Cells.Item(122, 11).Value = 3567
$ws.Cells.Item(122, 13).Value = -1117

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1683.3334
$ws.Cells.Item(126, 9).Value = 1550
$ws.Cells.Item(126, 11).Value = 4650
$ws.Cells.Item(126, 13).Value = -2180

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 910.7692
$ws.Cells.Item(136, 9).Value = 912.6667
$ws.Cells.Item(136, 11).Value = 2738.0001
$ws.Cells.Item(136, 13).Value = -188.0001000000002
